$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 11317
$wsExhibition.Range("F9").Value = 11255

# Sheet "全部类型" (sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 11317
$wsAll.Range("F11").Value = 11255
